$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" to H1, matching the style/format of the other header cells (G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the new "Save" data column values for the existing rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
